$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

Set-TextValue "D2" '42.002.07'
Set-TextValue "E2" '  -2.95%  '
Set-TextValue "D3" '2.234.77'
Set-TextValue "E3" '  -3.16%  '
Set-TextValue "E4" '  +0.04%  '
Set-TextValue "D5" '246.17'
Set-TextValue "E5" '  -2.71%  '
Set-TextValue "E6" '  -1.81%  '
Set-TextValue "D7" '75.81'
Set-TextValue "E7" '  +1.75%  '
Set-TextValue "E8" '  +0.11%  '
Set-TextValue "D9" '0.625'
Set-TextValue "E9" '  -4.40%  '
Set-TextValue "D10" '40.32'
Set-TextValue "E10" '  +1.95%  '
Set-TextValue "D11" '0.0951'
Set-TextValue "E11" '  -4.31%  '
Set-TextValue "D12" '7.12'
Set-TextValue "E12" '  -5.42%  '
Set-TextValue "E13" '  -2.16%  '
Set-TextValue "D14" '2.570.66'
Set-TextValue "E14" '  -3.11%  '
Set-TextValue "D15" '14.82'
Set-TextValue "E15" '  -5.08%  '
Set-TextValue "D16" '0.861'
Set-TextValue "E16" '  -2.77%  '
Set-TextValue "D17" '2.237.31'
Set-TextValue "E17" '  -3.22%  '
Set-TextValue "D18" '41.862.64'
Set-TextValue "E18" '  -2.91%  '
Set-TextValue "E19" '  -2.85%  '
Set-TextValue "E20" '  -3.46%  '
Set-TextValue "D21" '71.60'
Set-TextValue "E21" '  -1.77%  '
Set-TextValue "E22" '  +2.34%  '
Set-TextValue "D23" '231.10'
Set-TextValue "E23" '  -2.14%  '
Set-TextValue "D24" '11.44'
Set-TextValue "E24" '  -2.01%  '
Set-TextValue "E25" '  -0.01%  '
Set-TextValue "E26" '  -5.98%  '
Set-TextValue "E27" '  -5.70%  '
Set-TextValue "D28" '7.20'
Set-TextValue "E28" '  +12.16%  '
Set-TextValue "D29" '2.16'
Set-TextValue "E29" '  -1.75%  '
Set-TextValue "D30" '169.25'
Set-TextValue "E30" '  +0.89%  '
Set-TextValue "D31" '20.53'
Set-TextValue "E31" '  -3.16%  '
Set-TextValue "D32" '33.39'
Set-TextValue "E32" '  +3.73%  '
Set-TextValue "D33" '0.0838'
Set-TextValue "E33" '  +2.68%  '
Set-TextValue "E34" '  -5.88%  '
Set-TextValue "E35" '  -0.89%  '
Set-TextValue "D36" '4.51'
Set-TextValue "E36" '  -2.03%  '
Set-TextValue "E37" '  +1.19%  '
Set-TextValue "E38" '  -3.35%  '
Set-TextValue "E39" '  -8.06%  '
Set-TextValue "E40" '  -1.59%  '
Set-TextValue "E41" '  -7.98%  '
Set-TextValue "E42" '  +13.00%  '
Set-TextValue "E43" '  -6.10%  '
Set-TextValue "D44" '60.49'
Set-TextValue "E44" '  -3.27%  '
Set-TextValue "D45" '8.65'
Set-TextValue "E45" '  -5.76%  '
Set-TextValue "D46" '0.100'
Set-TextValue "E46" '  -3.74%  '
Set-TextValue "E48" '  -4.47%  '
Set-TextValue "D49" '4.36'
Set-TextValue "E49" '  -11.10%  '
Set-TextValue "E51" '  -2.77%  '
